$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting
# (values like '1.00' or '0.999' would otherwise be reinterpreted as numbers)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.896.59'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.497.13'
$ws.Range('E3').Value = '  -0.79%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.37'
$ws.Range('E5').Value = '  +0.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.13'
$ws.Range('E6').Value = '  -1.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.566'
$ws.Range('E8').Value = '  +0.75%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.522.70'
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').Value = '  +1.95%  '
$ws.Range('E11').Value = '  -0.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.36'
$ws.Range('E12').Value = '  -0.80%  '
$ws.Range('E13').Value = '  -2.08%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.973.33'
$ws.Range('E14').Value = '  +0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.21'
$ws.Range('E15').Value = '  +1.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.111.06'
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('E17').Value = '  -0.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.513.39'
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.06'
$ws.Range('E19').Value = '  +1.58%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.26'
$ws.Range('E20').Value = '  +0.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '325.17'
$ws.Range('E21').Value = '  +1.28%  '
$ws.Range('E22').Value = '  +0.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.86'
$ws.Range('E23').Value = '  +0.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '64.97'
$ws.Range('E24').Value = '  +3.43%  '
$ws.Range('E25').Value = '  -0.84%  '
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.01'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('E28').Value = '  -1.73%  '
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0774'
$ws.Range('E30').Value = '  +1.26%  '
$ws.Range('E31').Value = '  -1.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '167.80'
$ws.Range('E32').Value = '  +4.57%  '
$ws.Range('E33').Value = '  +3.83%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  -3.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.58'
$ws.Range('E36').Value = '  +0.66%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.12'
$ws.Range('E37').Value = '  -1.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.56'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.86'
$ws.Range('E39').Value = '  -0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.830'
$ws.Range('E40').Value = '  +3.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.63'
$ws.Range('E41').Value = '  +0.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.27'
$ws.Range('E42').Value = '  +0.43%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '280.89'
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.994'
$ws.Range('E44').Value = '  -0.39%  '
$ws.Range('E45').Value = '  +1.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.87'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '127.93'
$ws.Range('E47').Value = '  +4.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0930'
$ws.Range('E48').Value = '  +0.18%  '
$ws.Range('E49').Value = '  +0.85%  '
$ws.Range('E50').Value = '  +0.27%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.35'
$ws.Range('E51').Value = '  -0.14%  '

Write-Output "Update complete"
